$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.592.91"
$ws.Range("E2").Value = "  +2.79%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.912.63"
$ws.Range("E3").Value = "  +5.79%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.38"
$ws.Range("E5").Value = "  +1.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.0000"
$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5038"
$ws.Range("E7").Value = "  +2.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3960"
$ws.Range("E8").Value = "  +2.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09629"
$ws.Range("E9").Value = "  -1.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.157"
$ws.Range("E10").Value = "  +5.41%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.57"
$ws.Range("E11").Value = "  +2.68%  "

$ws.Range("E12").Value = "  +2.10%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.20"
$ws.Range("E13").Value = "  +3.81%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.904.26"
$ws.Range("E14").Value = "  +5.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.528"
$ws.Range("E15").Value = "  +3.45%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  -0.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001133"
$ws.Range("E17").Value = "  -0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.75"
$ws.Range("E18").Value = "  +1.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06641"
$ws.Range("E19").Value = "  +0.91%  "

$ws.Range("E20").Value = "  +5.31%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.289"
$ws.Range("E22").Value = "  +6.03%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.641.11"
$ws.Range("E23").Value = "  +2.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.39"
$ws.Range("E24").Value = "  +2.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("E25").Value = "  +1.82%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.755"
$ws.Range("E26").Value = "  +15.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.129.89"
$ws.Range("E27").Value = "  +5.52%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.43"
$ws.Range("E28").Value = "  +4.07%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "159.35"
$ws.Range("E29").Value = "  +1.56%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.46"
$ws.Range("E30").Value = "  +1.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.107"
$ws.Range("E31").Value = "  +6.55%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1072"
$ws.Range("E32").Value = "  +1.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.718"
$ws.Range("E33").Value = "  +2.91%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.622"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.761"
$ws.Range("E35").Value = "  +8.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06788"
$ws.Range("E36").Value = "  +0.88%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02447"
$ws.Range("E37").Value = "  +6.11%  "

$ws.Range("E38").Value = "  +4.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.082"
$ws.Range("E39").Value = "  +3.22%  "

$ws.Range("E40").Value = "  +3.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6412"
$ws.Range("E41").Value = "  +3.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.197"
$ws.Range("E42").Value = "  +5.12%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9994"
$ws.Range("E43").Value = "  -0.06%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.79"
$ws.Range("E44").Value = "  +5.98%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6058"
$ws.Range("E45").Value = "  +3.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.285"
$ws.Range("E46").Value = "  +0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.653"
$ws.Range("E47").Value = "  -1.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.035"
$ws.Range("E48").Value = "  +5.60%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.62"
$ws.Range("E49").Value = "  +1.81%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.210"
$ws.Range("E50").Value = "  +3.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06942"
$ws.Range("E51").Value = "  +2.29%  "

Write-Output "Update complete"